# "Added last minute updates"
# - paragraph 1: add a pBdr (top/left/bottom/right, 5-twip space, no line) and
#   bump the left indent from 120 -> 225 twips (6pt -> 11.25pt), matching the
#   formatting already used by the third paragraph in this doc.
# - paragraph 1: rename the placeholder id from
#   **ID__AFFARS_5301_topic_20__ID** to **ID__AFFARS_5301_602_1__ID**, and
#   drop the now-unneeded trailing-space run that followed it.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- paragraph formatting -------------------------------------------------
$p1.Format.Borders.DistanceFromTop    = 5
$p1.Format.Borders.DistanceFromLeft   = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight  = 5

$p1.Format.LeftIndent = 11.25   # 225 twips

# --- run content -----------------------------------------------------------
# Locate the placeholder id text within paragraph 1 (scoped search so we
# can't accidentally match anything outside this paragraph).
$idRange = $p1.Range.Duplicate
$null = $idRange.Find.Execute("**ID__AFFARS_5301_topic_20__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The trailing " " run sits right after the id run and before the paragraph
# mark; remove it first (it becomes dead weight once the id text changes).
$spaceRange = $d.Range($idRange.End, $p1.Range.End - 1)
$spaceRange.Delete()

# Now rewrite the id run's text in place.
$idRange.Text = "**ID__AFFARS_5301_602_1__ID**"
